$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(44344, 0, 2, 15.70475068708284),
    @(44345, 0, 2, 15.70475068708284),
    @(44346, 0, 2, 15.70475068708284),
    @(44347, 1, 1, 7.852375343541421),
    @(44348, 0, 1, 7.852375343541421),
    @(44349, 0, 1, 7.852375343541421),
    @(44350, 0, 1, 7.852375343541421),
    @(44351, 4, 5, 39.26187671770711),
    @(44352, 0, 5, 39.26187671770711),
    @(44353, 0, 5, 39.26187671770711),
    @(44354, 0, 4, 31.40950137416569),
    @(44355, 0, 4, 31.40950137416569),
    @(44356, 0, 4, 31.40950137416569),
    @(44357, 0, 4, 31.40950137416569),
    @(44358, 0, 0, 0),
    @(44359, 0, 0, 0),
    @(44360, 0, 0, 0),
    @(44361, 0, 0, 0),
    @(44362, 0, 0, 0),
    @(44363, 0, 0, 0),
    @(44364, 1, 1, 7.852375343541421),
    @(44365, 1, 2, 15.70475068708284),
    @(44366, 0, 2, 15.70475068708284),
    @(44367, 0, 2, 15.70475068708284),
    @(44368, 0, 2, 15.70475068708284),
    @(44369, 0, 2, 15.70475068708284),
    @(44370, 0, 2, 15.70475068708284),
    @(44371, 1, 2, 15.70475068708284),
    @(44372, 2, 3, 23.55712603062426),
    @(44373, 0, 3, 23.55712603062426),
    @(44374, 3, 6, 47.11425206124853),
    @(44375, 3, 9, 70.67137809187278)
)

$startRow = 270

# Copy the formatting of the last existing data row (269) down through
# the newly appended rows so column A keeps its date style (s="2").
$ws.Range("A269:D269").Copy()
$ws.Range("A270:D301").PasteSpecial(-4122)

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $startRow + $i
    $rowVals = $data[$i]
    $ws.Cells.Item($r, 1).Value = $rowVals[0]
    $ws.Cells.Item($r, 2).Value = $rowVals[1]
    $ws.Cells.Item($r, 3).Value = $rowVals[2]
    $ws.Cells.Item($r, 4).Value = $rowVals[3]
}
